# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced to
# Text ("@") before assignment so Excel keeps them as text instead of coercing
# them to numeric values (matches the inlineStr type in the source XML).
$textForceCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values
$ws.Range("D2").Value = '30.725.12'
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").Value = '1.884.97'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '239.65'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '0.4824'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = '0.2847'
$ws.Range("E8").Value = '  -1.43%  '
$ws.Range("D9").Value = '0.06553'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").Value = '1.991.88'
$ws.Range("E10").Value = '  +6.24%  '
$ws.Range("D11").Value = '0.07502'
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").Value = '16.73'
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("D13").Value = '5.138'
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").Value = '89.09'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '0.6697'
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = '30.672.80'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").Value = '13.41'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '0.9997'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007656'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '234.52'
$ws.Range("E20").Value = '  +18.49%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.181.10'
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").Value = '5.346'
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '6.239'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("D25").Value = '9.365'
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").Value = '168.65'
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("D27").Value = '18.72'
$ws.Range("E27").Value = '  +2.69%  '
$ws.Range("D28").Value = '1.963'
$ws.Range("E28").Value = '  +1.59%  '
$ws.Range("D29").Value = '1.449'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").Value = '0.09599'
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("D31").Value = '4.376'
$ws.Range("E31").Value = '  +2.82%  '
$ws.Range("D32").Value = '4.063'
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("D33").Value = '0.05065'
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("D34").Value = '1.222'
$ws.Range("E34").Value = '  +6.69%  '
$ws.Range("D35").Value = '0.7543'
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("D36").Value = '2.705'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '0.01873'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").Value = '2.632'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.112'
$ws.Range("E39").Value = '  +1.66%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.9196'
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").Value = '106.48'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '0.4314'
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").Value = '5.824'
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = '7.498'
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '65.17'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1298'
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("D48").Value = '1.491'
$ws.Range("E48").Value = '  -5.18%  '
$ws.Range("D49").Value = '8.990'
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("D50").Value = '34.11'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").Value = '0.3900'
$ws.Range("E51").Value = '  +0.68%  '

# Restore default (unstyled) cell style for the text-forced cells so the only
# change versus the original is the cell content, not its formatting.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
